$wb = $excel.ActiveWorkbook

# --- Sheet "ATS Accuracy" ---
$ws1 = $wb.Worksheets.Item("ATS Accuracy")

# Row 2: 5-fireball ATS
$ws1.Range("C2").Value = 84
$ws1.Range("D2").Value = 87
$ws1.Range("E2").Value = 96.59999999999999

# Row 3: 4-fireball ATS
$ws1.Range("C3").Value = 60
$ws1.Range("D3").Value = 63

# Row 4: 3-fireball ATS
$ws1.Range("B4").Value = 4
$ws1.Range("C4").Value = 15
$ws1.Range("D4").Value = 19
$ws1.Range("E4").Value = 78.90000000000001

# --- Sheet "Total Accuracy" ---
$ws2 = $wb.Worksheets.Item("Total Accuracy")

# Row 2: 5-fireball Total
$ws2.Range("C2").Value = 75
$ws2.Range("D2").Value = 80
$ws2.Range("E2").Value = 93.8

# Row 3: 4-fireball Total
$ws2.Range("B3").Value = 6
$ws2.Range("D3").Value = 67
$ws2.Range("E3").Value = 91
